$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The current last paragraph is "Finished a first version of the preview", which
# also carries the _GoBack bookmark markers at its end. The edit appends three new
# paragraphs after it and moves the bookmark markers to the end of the last of
# those new paragraphs. Re-write the existing last paragraph without the bookmark,
# then append the new paragraphs (the final one carrying the bookmark).

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$xmlFinished = "<w:p $wNs $w14Ns w14:paraId=`"28A48E33`" w14:textId=`"76781909`" w:rsidR=`"00CF05D1`" w:rsidRPr=`"00012996`" w:rsidRDefault=`"00CF05D1`" w:rsidP=`"00CF05D1`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"9`"/></w:numPr></w:pPr><w:r><w:t>Finished a first version of the preview</w:t></w:r></w:p>"
$lastRange.InsertXML($xmlFinished)

# Append "9/11/18" as a plain paragraph (no style / list numbering).
$end = $d.Content
$end.Collapse(0)
$xmlDate = "<w:p $wNs><w:r><w:t>9/11/18</w:t></w:r></w:p>"
$end.InsertXML($xmlDate)

# Append the "Removed meta files..." list item, including the rendered-page-break
# marker that was present in the target revision.
$end = $d.Content
$end.Collapse(0)
$xmlRemoved = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"9`"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Removed meta files for unity from the git ignore.</w:t></w:r></w:p>"
$end.InsertXML($xmlRemoved)

# Append the "Updated nodes..." list item (split across three runs, matching the
# source revision), carrying the relocated _GoBack bookmark markers at its end.
$end = $d.Content
$end.Collapse(0)
$xmlUpdated = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"9`"/></w:numPr></w:pPr><w:r><w:t>U</w:t></w:r><w:r><w:t>pdated</w:t></w:r><w:r><w:t xml:space=`"preserve`"> nodes so that they return running rather than continuing</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$end.InsertXML($xmlUpdated)
